$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update the phone number for Eziechiele (row 3, column C)
$ws.Range("C3").Value = 6923378500

# Delete rows 4 through 11 (the extra generated names) so only the
# header + 2 data rows remain
$ws.Range("A4:C11").EntireRow.Delete()

# Update the selection to match the authored state
$ws.Range("C3").Select()
